# Update the "Clasificación" league table with the results of a newly
# played round: each player's games played (PJ) goes from 2 to 3, and the
# related stats (wins/draws/losses, goals, points, etc.) are updated
# accordingly.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Clasificación")

# Row 2 - David
$ws.Range("C2").Value = 3.0
$ws.Range("D2").Value = 2.0
$ws.Range("J2").Value = 2.0
$ws.Range("M2").Value = 6.0

# Row 3 - Pedro
$ws.Range("C3").Value = 3.0
$ws.Range("F3").Value = 3.0
$ws.Range("H3").Value = 6.0
$ws.Range("I3").Value = -5.0

# Row 4 - Adonay
$ws.Range("C4").Value = 3.0
$ws.Range("D4").Value = 3.0
$ws.Range("G4").Value = 4.0
$ws.Range("I4").Value = 4.0
$ws.Range("J4").Value = 2.0
$ws.Range("M4").Value = 10.0

# Row 5 - Richard
$ws.Range("C5").Value = 3.0
$ws.Range("F5").Value = 2.0
$ws.Range("H5").Value = 2.0
$ws.Range("I5").Value = -2.0

# Row 6 - Iván
$ws.Range("C6").Value = 3.0
$ws.Range("D6").Value = 2.0
$ws.Range("G6").Value = 6.0
$ws.Range("I6").Value = 4.0
$ws.Range("L6").Value = 2.0
$ws.Range("M6").Value = 10.0

# Row 7 - Nico
$ws.Range("C7").Value = 3.0
$ws.Range("F7").Value = 1.0
$ws.Range("H7").Value = 1.0

# Row 8 - Nicolás
$ws.Range("C8").Value = 3.0
$ws.Range("F8").Value = 2.0
$ws.Range("H8").Value = 4.0
$ws.Range("I8").Value = -3.0
$ws.Range("M8").Value = 4.0

# Row 9 - Vicente
$ws.Range("C9").Value = 3.0
$ws.Range("D9").Value = 2.0
$ws.Range("G9").Value = 2.0
$ws.Range("I9").Value = 1.0
$ws.Range("J9").Value = 2.0
$ws.Range("M9").Value = 6.0
